# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes,
# plus the FraxShare/PaxDollar row swap (rows 41-42: B/C/D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the value to be written as text, even when it looks numeric
    # (e.g. "19.30", "1.00"), then restore the default "Normal" style so
    # the cell formatting stays identical to the original workbook.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.259.34'
Set-TextValue $ws.Range("E2") '  +0.30%  '

Set-TextValue $ws.Range("D3") '1.590.16'
Set-TextValue $ws.Range("E3") '  +0.59%  '

Set-TextValue $ws.Range("E4") '  -0.21%  '

Set-TextValue $ws.Range("D5") '212.69'
Set-TextValue $ws.Range("E5") '  +1.53%  '

Set-TextValue $ws.Range("E6") '  +1.01%  '

Set-TextValue $ws.Range("E7") '  -0.17%  '

Set-TextValue $ws.Range("E8") '  +0.20%  '

Set-TextValue $ws.Range("E9") '  -0.17%  '

Set-TextValue $ws.Range("D10") '19.30'
Set-TextValue $ws.Range("E10") '  -1.11%  '

Set-TextValue $ws.Range("D11") '0.0848'
Set-TextValue $ws.Range("E11") '  +0.38%  '

Set-TextValue $ws.Range("D12") '1.812.60'
Set-TextValue $ws.Range("E12") '  +0.55%  '

Set-TextValue $ws.Range("D13") '1.579.53'
Set-TextValue $ws.Range("E13") '  -0.14%  '

Set-TextValue $ws.Range("E14") '  -0.09%  '

Set-TextValue $ws.Range("D15") '0.521'
Set-TextValue $ws.Range("E15") '  +1.27%  '

Set-TextValue $ws.Range("D16") '64.40'
Set-TextValue $ws.Range("E16") '  -0.14%  '

Set-TextValue $ws.Range("D17") '26.257.22'
Set-TextValue $ws.Range("E17") '  +0.25%  '

Set-TextValue $ws.Range("E18") '  -1.01%  '

Set-TextValue $ws.Range("D19") '7.43'
Set-TextValue $ws.Range("E19") '  +2.21%  '

Set-TextValue $ws.Range("D20") '213.24'
Set-TextValue $ws.Range("E20") '  +2.86%  '

Set-TextValue $ws.Range("E21") '  -0.16%  '

Set-TextValue $ws.Range("E22") '  +0.59%  '

Set-TextValue $ws.Range("E23") '  +1.17%  '

Set-TextValue $ws.Range("E24") '  -2.55%  '

Set-TextValue $ws.Range("D25") '144.64'
Set-TextValue $ws.Range("E25") '  +0.04%  '

Set-TextValue $ws.Range("E26") '  -0.29%  '

Set-TextValue $ws.Range("E27") '  +0.78%  '

Set-TextValue $ws.Range("E28") '  -0.66%  '

Set-TextValue $ws.Range("E29") '  -0.28%  '

Set-TextValue $ws.Range("E30") '  -1.24%  '

Set-TextValue $ws.Range("E31") '  +0.85%  '

Set-TextValue $ws.Range("D32") '3.21'
Set-TextValue $ws.Range("E32") '  -0.16%  '

Set-TextValue $ws.Range("D33") '2.95'
Set-TextValue $ws.Range("E33") '  -0.23%  '

Set-TextValue $ws.Range("D34") '1.335.81'
Set-TextValue $ws.Range("E34") '  +4.68%  '

Set-TextValue $ws.Range("D35") '2.44'
Set-TextValue $ws.Range("E35") '  -0.94%  '

Set-TextValue $ws.Range("E36") '  -0.74%  '

Set-TextValue $ws.Range("E37") '  -2.73%  '

Set-TextValue $ws.Range("E38") '  +0.08%  '

Set-TextValue $ws.Range("D39") '0.816'
Set-TextValue $ws.Range("E39") '  -0.12%  '

Set-TextValue $ws.Range("D40") '1.04'
Set-TextValue $ws.Range("E40") '  -3.17%  '

Set-TextValue $ws.Range("B41") 'PaxDollar'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D41") '1.00'
Set-TextValue $ws.Range("E41") '  -0.18%  '

Set-TextValue $ws.Range("B42") 'FraxShare'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D42") '5.70'
Set-TextValue $ws.Range("E42") '  +3.15%  '

Set-TextValue $ws.Range("D43") '2.13'
Set-TextValue $ws.Range("E43") '  -0.15%  '

Set-TextValue $ws.Range("D44") '0.763'
Set-TextValue $ws.Range("E44") '  -0.09%  '

Set-TextValue $ws.Range("D45") '1.724.90'
Set-TextValue $ws.Range("E45") '  +0.46%  '

Set-TextValue $ws.Range("D46") '61.76'
Set-TextValue $ws.Range("E46") '  -0.95%  '

Set-TextValue $ws.Range("D47") '86.59'
Set-TextValue $ws.Range("E47") '  -2.66%  '

Set-TextValue $ws.Range("E49") '  -0.69%  '

Set-TextValue $ws.Range("D50") '0.0975'
Set-TextValue $ws.Range("E50") '  -2.88%  '

Set-TextValue $ws.Range("D51") '1.00'
Set-TextValue $ws.Range("E51") '  -0.36%  '
